$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://www.jonglaan.nl/"
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = "[('geolocation', ''), ('midi', ''), ('notifications', ''), ('push', ''), ('sync-xhr', ''), ('microphone', ''), ('camera', ''), ('magnetometer', ''), ('gyroscope', ''), ('speaker', 'self'), ('vibrate', ''), ('fullscreen', 'selfhttps://*.youtube.comhttps://*.youtube-nocookie.comhttps://*.youtu.be'), ('payment', '')]"
$ws.Range("E2").Value = "[['microphone *', 'about:blank']]"
$ws.Range("F2").Value = "[]"
$ws.Range("J2").Value = "['about:blank']"
